$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new rows above row 384. This shifts the existing rows
# 384:496 down to 387:499 (values + formatting move together, matching
# Excel's native "Insert Copied/Blank Rows" behaviour).
$ws.Rows("384:386").Insert()

# Populate the three freshly-inserted rows with the new (most recent)
# weekly price report for "Cuatro cascos verde".
$common = @{
    A = 2
    B = "Comercializadora del Agro de Limarí"
    C = "Coquimbo"
    D = 44468
    E = 4
    F = 100112002
    G = "Pimiento"
    H = "Cuatro cascos verde"
    N = "`$/caja 18 kilos"
    O = "Provincia de Limarí"
    Q = 18
    R = "Hortaliza"
}

$rows = @(
    @{ Row = 384; I = "Primera"; J = 800; K = 32000; L = 33000; M = 32500; P = 1806 },
    @{ Row = 385; I = "Segunda"; J = 800; K = 29000; L = 30000; M = 29500; P = 1639 },
    @{ Row = 386; I = "Tercera"; J = 500; K = 25000; L = 26000; M = 25500; P = 1417 }
)

foreach ($r in $rows) {
    $rowNum = $r.Row
    $ws.Cells.Item($rowNum, 1).Value = $common.A
    $ws.Cells.Item($rowNum, 2).Value = $common.B
    $ws.Cells.Item($rowNum, 3).Value = $common.C
    $ws.Cells.Item($rowNum, 4).Value = $common.D
    $ws.Cells.Item($rowNum, 5).Value = $common.E
    $ws.Cells.Item($rowNum, 6).Value = $common.F
    $ws.Cells.Item($rowNum, 7).Value = $common.G
    $ws.Cells.Item($rowNum, 8).Value = $common.H
    $ws.Cells.Item($rowNum, 9).Value = $r.I
    $ws.Cells.Item($rowNum, 10).Value = $r.J
    $ws.Cells.Item($rowNum, 11).Value = $r.K
    $ws.Cells.Item($rowNum, 12).Value = $r.L
    $ws.Cells.Item($rowNum, 13).Value = $r.M
    $ws.Cells.Item($rowNum, 14).Value = $common.N
    $ws.Cells.Item($rowNum, 15).Value = $common.O
    $ws.Cells.Item($rowNum, 16).Value = $r.P
    $ws.Cells.Item($rowNum, 17).Value = $common.Q
    $ws.Cells.Item($rowNum, 18).Value = $common.R
}
